$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix two values in row 3 (small odds corrections)
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63

# Add new row 4 data
$ws.Range("A4").Value = "fakytOLi"
$ws.Range("B4").Value = "15/10/2024"
$ws.Range("C4").Value = "15:45"
$ws.Range("D4").Value = "WALES - CYMRU PREMIER"
$ws.Range("E4").Value = "Cardiff Metropolitan"
$ws.Range("F4").Value = "Penybont"
$ws.Range("G4").Value = 3.95
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 1.85
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.47
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 6.7
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 2.92
$ws.Range("Q4").Value = 2.07
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 1.42
$ws.Range("T4").Value = 2.65
$ws.Range("U4").Value = 1.9
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 10.5
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 13.5
$ws.Range("Z4").Value = 65
$ws.Range("AA4").Value = 40
$ws.Range("AB4").Value = 45
$ws.Range("AC4").Value = 6.7
$ws.Range("AD4").Value = 6.4
$ws.Range("AE4").Value = 16
$ws.Range("AF4").Value = 80
$ws.Range("AG4").Value = 6.3
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 8.5
$ws.Range("AJ4").Value = 15.5
$ws.Range("AK4").Value = 16
$ws.Range("AL4").Value = 30
$ws.Range("AM4").Value = 700
$ws.Range("AN4").Value = 5.8
$ws.Range("AO4").Value = 23
$ws.Range("AP4").Value = 30
$ws.Range("AQ4").Value = 120
$ws.Range("AR4").Value = 175
$ws.Range("AS4").Value = 400
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 7.4
$ws.Range("AV4").Value = 70
$ws.Range("AW4").Value = 3.7
$ws.Range("AX4").Value = 9.5
$ws.Range("AY4").Value = 19.5
$ws.Range("AZ4").Value = 35
$ws.Range("BA4").Value = 75
$ws.Range("BB4").Value = 300
$ws.Range("BC4").Value = 51
$ws.Range("BD4").Value = 51
